$wb = $excel.ActiveWorkbook

# New PC-1 / PC-2 values for rows 2-6, column B (PC-1) and C (PC-2), per sheet.
$newValues = @{
    1 = @{
        2 = @(0.5115821815672132, 0.02175819009777252)
        3 = @(0.4054111960378419, -0.7118471303185286)
        4 = @(0.4113086412639908, 0.1341045258859178)
        5 = @(0.3901037336585669, 0.6834080668637618)
        6 = @(0.5025631425376446, -0.08814553025193195)
    }
    2 = @{
        2 = @(-0.5186840863940005, -0.01520805295325373)
        3 = @(-0.2776044042943789, -0.8610860848542273)
        4 = @(-0.4606501046128745, 0.1266835982382657)
        5 = @(-0.4213634355599578, 0.4909035238897989)
        6 = @(-0.5139620117588853, -0.03555932740888018)
    }
    3 = @{
        2 = @(-0.5194204750270668, -0.07806132019184023)
        3 = @(-0.3448150747342947, 0.6989225872624897)
        4 = @(-0.4622486903490475, -0.3276280529488359)
        5 = @(-0.4139741641828348, -0.5233253765996744)
        6 = @(-0.4756642450493342, 0.3524259587941533)
    }
    4 = @{
        2 = @(-0.5660876030972491, 0.1083208302699449)
        3 = @(-0.1928640120986571, 0.6946371660720148)
        4 = @(-0.4930486747734305, -0.3311659381744667)
        5 = @(-0.3827057494038196, -0.5298331396828329)
        6 = @(-0.5027798843770949, 0.3396347607768182)
    }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Remove the last data row (row 7, "mean_length_depth_width") entirely.
    $ws.Rows.Item(7).Delete()

    # Update the remaining PC-1 / PC-2 values.
    $sheetValues = $newValues[$i]
    foreach ($r in $sheetValues.Keys) {
        $pair = $sheetValues[$r]
        $ws.Cells.Item($r, 2).Value = $pair[0]
        $ws.Cells.Item($r, 3).Value = $pair[1]
    }
}
